$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 35703
$ws.Range("D2").Value = 51692268
$ws.Range("C3").Value = 87031
$ws.Range("D3").Value = 127709223
$ws.Range("C4").Value = 29842
$ws.Range("D4").Value = 44227845
$ws.Range("C5").Value = 8249
$ws.Range("D5").Value = 12265510
$ws.Range("C6").Value = 1775
$ws.Range("D6").Value = 2639106
$ws.Range("C7").Value = 133
$ws.Range("D7").Value = 194593
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 13500
$ws.Range("C11").Value = 39087
$ws.Range("D11").Value = 53156134
$ws.Range("C12").Value = 9193
$ws.Range("D12").Value = 13305791
$ws.Range("C13").Value = 24985
$ws.Range("D13").Value = 36666856
$ws.Range("C14").Value = 7992
$ws.Range("D14").Value = 11868031
$ws.Range("C15").Value = 2046
$ws.Range("D15").Value = 3044150
$ws.Range("C16").Value = 378
$ws.Range("D16").Value = 556123
$ws.Range("C19").Value = 9699
$ws.Range("D19").Value = 12883443
$ws.Range("C20").Value = 12828
$ws.Range("D20").Value = 18537638
$ws.Range("C21").Value = 30561
$ws.Range("D21").Value = 44894843
$ws.Range("C22").Value = 9899
$ws.Range("D22").Value = 14721695
$ws.Range("C23").Value = 2523
$ws.Range("D23").Value = 3754763
$ws.Range("C24").Value = 454
$ws.Range("D24").Value = 675345
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 45953
$ws.Range("C26").Value = 11154
$ws.Range("D26").Value = 14954265
$ws.Range("C27").Value = 7311
$ws.Range("D27").Value = 10595197
$ws.Range("C28").Value = 21666
$ws.Range("D28").Value = 31811900
$ws.Range("C29").Value = 7521
$ws.Range("D29").Value = 11196209
$ws.Range("C30").Value = 1871
$ws.Range("D30").Value = 2792144
$ws.Range("C31").Value = 316
$ws.Range("D31").Value = 471415
$ws.Range("C32").Value = 27
$ws.Range("D32").Value = 40393
$ws.Range("C33").Value = 7928
$ws.Range("D33").Value = 10503967
$ws.Range("C34").Value = 2925
$ws.Range("D34").Value = 4216484
$ws.Range("C35").Value = 7259
$ws.Range("D35").Value = 10610777
$ws.Range("C36").Value = 2930
$ws.Range("D36").Value = 4339288
$ws.Range("C37").Value = 779
$ws.Range("D37").Value = 1161263
$ws.Range("C38").Value = 140
$ws.Range("D38").Value = 208232
$ws.Range("C40").Value = 2243
$ws.Range("D40").Value = 3025971
$ws.Range("C41").Value = 16408
$ws.Range("D41").Value = 23744106
$ws.Range("C42").Value = 49034
$ws.Range("D42").Value = 71926435
$ws.Range("C43").Value = 18318
$ws.Range("D43").Value = 27214255
$ws.Range("C44").Value = 5333
$ws.Range("D44").Value = 7945298
$ws.Range("C45").Value = 1073
$ws.Range("D45").Value = 1600292
$ws.Range("C46").Value = 56
$ws.Range("D46").Value = 82348
$ws.Range("C49").Value = 15875
$ws.Range("D49").Value = 21189338
$ws.Range("C50").Value = 1750
$ws.Range("D50").Value = 2542379
$ws.Range("C51").Value = 6226
$ws.Range("D51").Value = 9162516
$ws.Range("C52").Value = 2161
$ws.Range("D52").Value = 3228250
$ws.Range("C53").Value = 705
$ws.Range("D53").Value = 1052805
$ws.Range("C54").Value = 163
$ws.Range("D54").Value = 242111
$ws.Range("C55").Value = 14
$ws.Range("D55").Value = 21000
$ws.Range("C56").Value = 5731
$ws.Range("D56").Value = 7912513
$ws.Range("C57").Value = 735
$ws.Range("D57").Value = 1077990
$ws.Range("C58").Value = 1868
$ws.Range("D58").Value = 2770126
$ws.Range("C59").Value = 756
$ws.Range("D59").Value = 1126453
$ws.Range("C60").Value = 256
$ws.Range("D60").Value = 383758
$ws.Range("C61").Value = 56
$ws.Range("D61").Value = 84000
$ws.Range("C63").Value = 1092
$ws.Range("D63").Value = 1545952
$ws.Range("C64").Value = 14634
$ws.Range("D64").Value = 21157056
$ws.Range("C65").Value = 43100
$ws.Range("D65").Value = 63113701
$ws.Range("C66").Value = 15174
$ws.Range("D66").Value = 22559656
$ws.Range("C67").Value = 4387
$ws.Range("D67").Value = 6534716
$ws.Range("C68").Value = 860
$ws.Range("D68").Value = 1279596
$ws.Range("C69").Value = 69
$ws.Range("D69").Value = 102189
$ws.Range("C71").Value = 14469
$ws.Range("D71").Value = 19137478
$ws.Range("C72").Value = 47557
$ws.Range("D72").Value = 69247673
$ws.Range("C73").Value = 136756
$ws.Range("D73").Value = 201594025
$ws.Range("C74").Value = 59841
$ws.Range("D74").Value = 89200723
$ws.Range("C75").Value = 18978
$ws.Range("D75").Value = 28361380
$ws.Range("C76").Value = 4284
$ws.Range("D76").Value = 6402620
$ws.Range("C77").Value = 229
$ws.Range("D77").Value = 338670
$ws.Range("C83").Value = 47068
$ws.Range("D83").Value = 64343779
$ws.Range("C84").Value = 4273
$ws.Range("D84").Value = 6196953
$ws.Range("C85").Value = 10924
$ws.Range("D85").Value = 16056465
$ws.Range("C86").Value = 3722
$ws.Range("D86").Value = 5547187
$ws.Range("C87").Value = 1297
$ws.Range("D87").Value = 1937989
$ws.Range("C88").Value = 271
$ws.Range("D88").Value = 404012
$ws.Range("C91").Value = 4945
$ws.Range("D91").Value = 6667045
$ws.Range("C92").Value = 1444
$ws.Range("D92").Value = 2087992
$ws.Range("C93").Value = 4771
$ws.Range("D93").Value = 7029734
$ws.Range("C94").Value = 1822
$ws.Range("D94").Value = 2715870
$ws.Range("C95").Value = 646
$ws.Range("D95").Value = 968141
$ws.Range("C96").Value = 163
$ws.Range("D96").Value = 243613
$ws.Range("C99").Value = 3174
$ws.Range("D99").Value = 4208909
$ws.Range("C100").Value = 536
$ws.Range("D100").Value = 799464
$ws.Range("C101").Value = 290
$ws.Range("D101").Value = 433165
$ws.Range("C102").Value = 100
$ws.Range("D102").Value = 150000
$ws.Range("C103").Value = 39
$ws.Range("D103").Value = 58500
$ws.Range("C104").Value = 19
$ws.Range("D104").Value = 28500
$ws.Range("C105").Value = 10293
$ws.Range("D105").Value = 14952645
$ws.Range("C106").Value = 28224
$ws.Range("D106").Value = 41486413
$ws.Range("C107").Value = 9454
$ws.Range("D107").Value = 14060120
$ws.Range("C108").Value = 2571
$ws.Range("D108").Value = 3833749
$ws.Range("C109").Value = 444
$ws.Range("D109").Value = 663482
$ws.Range("C110").Value = 44
$ws.Range("D110").Value = 66000
$ws.Range("C112").Value = 9353
$ws.Range("D112").Value = 12388428
$ws.Range("C113").Value = 28964
$ws.Range("D113").Value = 41806960
$ws.Range("C114").Value = 63726
$ws.Range("D114").Value = 93323933
$ws.Range("C115").Value = 20630
$ws.Range("D115").Value = 30668815
$ws.Range("C116").Value = 5771
$ws.Range("D116").Value = 8599213
$ws.Range("C117").Value = 1039
$ws.Range("D117").Value = 1554493
$ws.Range("C118").Value = 63
$ws.Range("D118").Value = 91920
$ws.Range("C121").Value = 24594
$ws.Range("D121").Value = 32921463
$ws.Range("C122").Value = 33815
$ws.Range("D122").Value = 48866400
$ws.Range("C123").Value = 73222
$ws.Range("D123").Value = 107171660
$ws.Range("C124").Value = 22797
$ws.Range("D124").Value = 33849205
$ws.Range("C125").Value = 6047
$ws.Range("D125").Value = 8989705
$ws.Range("C126").Value = 1121
$ws.Range("D126").Value = 1668146
$ws.Range("C127").Value = 53
$ws.Range("D127").Value = 77728
$ws.Range("C130").Value = 29852
$ws.Range("D130").Value = 39747675
$ws.Range("C131").Value = 12653
$ws.Range("D131").Value = 18328000
$ws.Range("C132").Value = 31195
$ws.Range("D132").Value = 45848159
$ws.Range("C133").Value = 11093
$ws.Range("D133").Value = 16483822
$ws.Range("C134").Value = 2816
$ws.Range("D134").Value = 4199620
$ws.Range("C135").Value = 451
$ws.Range("D135").Value = 670490
$ws.Range("C138").Value = 10359
$ws.Range("D138").Value = 13865316
$ws.Range("C139").Value = 33053
$ws.Range("D139").Value = 47775693
$ws.Range("C140").Value = 77510
$ws.Range("D140").Value = 113626017
$ws.Range("C141").Value = 23326
$ws.Range("D141").Value = 34681831
$ws.Range("C142").Value = 6051
$ws.Range("D142").Value = 9031499
$ws.Range("C143").Value = 1325
$ws.Range("D143").Value = 1974186
$ws.Range("C144").Value = 70
$ws.Range("D144").Value = 104630
$ws.Range("C146").Value = 27675
$ws.Range("D146").Value = 37467315
